$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 965.21875
$ws.Range("J17").Value = 948.79364
$ws.Range("L17").Value = 2846.38092
$ws.Range("N17").Value = -3182.38092
$ws.Range("H33").Value = 1371.5
$ws.Range("I33").Value = 1595.8889
$ws.Range("J33").Value = 698.3333
$ws.Range("K33").Value = 1595.8889
$ws.Range("L33").Value = 698.3333
$ws.Range("M33").Value = -1366.8889
$ws.Range("N33").Value = -1156.3333
$ws.Range("H53").Value = 199.75
$ws.Range("I53").Value = 200
$ws.Range("K53").Value = 200
$ws.Range("M53").Value = 437
$ws.Range("H62").Value = 4735
$ws.Range("I62").Value = 3557.5
$ws.Range("J62").Value = 6501.25
$ws.Range("K62").Value = 3557.5
$ws.Range("L62").Value = 6501.25
$ws.Range("M62").Value = -2933.5
$ws.Range("N62").Value = -7749.25
$ws.Range("H65").Value = 4735
$ws.Range("I65").Value = 3557.5
$ws.Range("J65").Value = 6501.25
$ws.Range("K65").Value = 17787.5
$ws.Range("L65").Value = 32506.25
$ws.Range("M65").Value = -14667.5
$ws.Range("N65").Value = -38746.25
$ws.Range("H116").Value = 7713.6665
$ws.Range("J116").Value = 4553
$ws.Range("L116").Value = 4553
$ws.Range("N116").Value = -11437
$ws.Range("H135").Value = 706.2564
$ws.Range("I135").Value = 545.97144
$ws.Range("J135").Value = 2108.75
$ws.Range("K135").Value = 4913.74296
$ws.Range("L135").Value = 18978.75
$ws.Range("M135").Value = -2378.74296
$ws.Range("N135").Value = -24048.75
$ws.Range("H138").Value = 14287462
$ws.Range("I138").Value = 25642336
$ws.Range("J138").Value = 2296.0967
$ws.Range("K138").Value = 76927008
$ws.Range("L138").Value = 6888.2901
$ws.Range("M138").Value = -76921868
$ws.Range("N138").Value = -17168.2901

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 333.44
$ws.Range("I32").Value = 292.8737
$ws.Range("J32").Value = 1104.2
$ws.Range("K32").Value = 292.8737
$ws.Range("L32").Value = 1104.2
$ws.Range("M32").Value = -5.873699999999985
$ws.Range("N32").Value = -1678.2
$ws.Range("H45").Value = 11846.2
$ws.Range("I45").Value = 20457
$ws.Range("K45").Value = 20457
$ws.Range("M45").Value = -20080

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 13290.333
$ws.Range("J103").Value = 13290.333
$ws.Range("L103").Value = 13290.333
$ws.Range("N103").Value = -15634.333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H22").Value = 4330.8335
$ws.Range("I22").Value = 4800
$ws.Range("J22").Value = 3861.6667
$ws.Range("K22").Value = 4800
$ws.Range("L22").Value = 3861.6667
$ws.Range("M22").Value = -4450
$ws.Range("N22").Value = -4561.6667
$ws.Range("H86").Value = 13361136
$ws.Range("I86").Value = 22245560
$ws.Range("J86").Value = 34499.5
$ws.Range("K86").Value = 22245560
$ws.Range("L86").Value = 34499.5
$ws.Range("M86").Value = -22244437
$ws.Range("N86").Value = -36745.5
$ws.Range("H89").Value = 13361136
$ws.Range("I89").Value = 22245560
$ws.Range("J89").Value = 34499.5
$ws.Range("K89").Value = 111227800
$ws.Range("L89").Value = 172497.5
$ws.Range("M89").Value = -111222184
$ws.Range("N89").Value = -183729.5
$ws.Range("H99").Value = 6944.1343
$ws.Range("I99").Value = 9153.154
$ws.Range("J99").Value = 6412.3335
$ws.Range("K99").Value = 9153.154
$ws.Range("L99").Value = 6412.3335
$ws.Range("M99").Value = -7655.154
$ws.Range("N99").Value = -9408.333500000001
$ws.Range("H107").Value = 1323.4286
$ws.Range("I107").Value = 1003.2
$ws.Range("J107").Value = 2124
$ws.Range("K107").Value = 1003.2
$ws.Range("L107").Value = 2124
$ws.Range("M107").Value = 916.8
$ws.Range("N107").Value = -5964
$ws.Range("H126").Value = 6944.1343
$ws.Range("I126").Value = 9153.154
$ws.Range("J126").Value = 6412.3335
$ws.Range("K126").Value = 27459.462
$ws.Range("L126").Value = 19237.0005
$ws.Range("M126").Value = -24989.462
$ws.Range("N126").Value = -24177.0005
$ws.Range("H132").Value = 7150.1113
$ws.Range("I132").Value = 3065.4092
$ws.Range("J132").Value = 25122.8
$ws.Range("K132").Value = 9196.2276
$ws.Range("L132").Value = 75368.39999999999
$ws.Range("M132").Value = -6666.2276
$ws.Range("N132").Value = -80428.39999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 419.72726
$ws.Range("I14").Value = 419.72726
$ws.Range("K14").Value = 1259.18178
$ws.Range("M14").Value = -1086.18178
$ws.Range("H33").Value = 241.16667
$ws.Range("J33").Value = 270.8
$ws.Range("L33").Value = 1624.8
$ws.Range("N33").Value = -2190.8
$ws.Range("H46").Value = 561.5714
$ws.Range("I46").Value = 449.33334
$ws.Range("J46").Value = 592.1818
$ws.Range("K46").Value = 1348.00002
$ws.Range("L46").Value = 1776.5454
$ws.Range("M46").Value = -1257.00002
$ws.Range("N46").Value = -1958.5454
$ws.Range("H81").Value = 1361.8334
$ws.Range("I81").Value = 1263.8889
$ws.Range("K81").Value = 3791.6667
$ws.Range("M81").Value = -2668.6667
$ws.Range("H84").Value = 1361.8334
$ws.Range("I84").Value = 1263.8889
$ws.Range("K84").Value = 11375.0001
$ws.Range("M84").Value = -5759.000099999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H101").Value = 20000
$ws.Range("J101").Value = 20000
$ws.Range("L101").Value = 20000
$ws.Range("N101").Value = -26490

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 63340.8
$ws.Range("J101").Value = 63340.8
$ws.Range("L101").Value = 63340.8
$ws.Range("N101").Value = -69830.8
$ws.Range("H126").Value = 10075.429
$ws.Range("I126").Value = 5897.3335
$ws.Range("J126").Value = 17596
$ws.Range("K126").Value = 17692.0005
$ws.Range("L126").Value = 52788
$ws.Range("M126").Value = -15222.0005
$ws.Range("N126").Value = -57728
$ws.Range("H132").Value = 3585.7632
$ws.Range("I132").Value = 2407.4285
$ws.Range("K132").Value = 7222.2855
$ws.Range("M132").Value = -4692.2855
$ws.Range("H136").Value = 2114.8108
$ws.Range("I136").Value = 2021.7428
$ws.Range("J136").Value = 3743.5
$ws.Range("K136").Value = 6065.2284
$ws.Range("L136").Value = 11230.5
$ws.Range("M136").Value = -3515.2284
$ws.Range("N136").Value = -16330.5
$ws.Range("H137").Value = 117506.86
$ws.Range("J137").Value = 117506.86
$ws.Range("L137").Value = 117506.86
$ws.Range("N137").Value = -127706.86
